# Apply the authored edit: rename/delete worksheets, update the
# "Correction Factor Calcs" note text, and add the per-county header row
# to the newly-named "Correction Transposed" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Sheet renames / deletion -------------------------------------------------
$wsTotal = $wb.Worksheets.Item("Total Workers_21")
$wsTotal.Name = "Total Worker"

$wsCalc = $wb.Worksheets.Item("Households by No. of Worker_21")
$wsCalc.Name = "Correction Factor Calcs"

$wsOldCorrection = $wb.Worksheets.Item("HHs by Workers Correction_21")
$wsOldCorrection.Delete()

$wsTransposed = $wb.Worksheets.Item("Households by No. of Worker (2)")
$wsTransposed.Name = "Correction Transposed"

# --- 2. Updated note text on "Correction Factor Calcs" ---------------------------
$wsCalc.Range("K5").Value = "ACS 2017-2021`nValues Inflated to Match Total Worker Table (B23025) from PUMS 2017-2021"

# --- 3. New per-county header row on "Correction Transposed" ---------------------
$wsTransposed.Range("F1").Value = "San Francisco"
$wsTransposed.Range("G1").Value = "San Mateo"
$wsTransposed.Range("H1").Value = "Santa Clara"
$wsTransposed.Range("I1").Value = "Alameda"
$wsTransposed.Range("J1").Value = "Contra Costa"
$wsTransposed.Range("K1").Value = "Solano"
$wsTransposed.Range("L1").Value = "Napa"
$wsTransposed.Range("M1").Value = "Sonoma"
$wsTransposed.Range("N1").Value = "Marin"

# --- 4. Selections / active sheet, matching the saved view state -----------------
$wsTotal.Range("B16").Select()

$wsCalc.Activate()
$wsCalc.Range("C13").Select()
$wsCalc.Range("I16").Select()

$wsTransposed.Range("E1:N5").Select()

$wsCalc.Activate()
